$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (D4: 2 -> 3, D6: 3 -> 2) ---
$ws.Range("D4").Value = 3
$ws.Range("D6").Value = 2

# --- New column K: "Protokoll" ---
# Header K3, styled like the other yellow header cells (copy format from J3)
$ws.Range("K3").Value = "Protokoll"
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# K4 = 1, styled like C4 (blue fill, centered) - value "1" style
$ws.Range("K4").Value = 1
$ws.Range("C4").Copy()
$ws.Range("K4").PasteSpecial(-4122)

# K5 = 3, new style: solid white (theme Background 1) fill, centered, bordered.
# Start from a solid-filled centered cell (C4) so Pattern/Border/Alignment/Font
# carry over, then recolor the fill to the theme "Background 1" color.
$ws.Range("K5").Value = 3
$ws.Range("C4").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Interior.ThemeColor = 2

# K6 = 2, styled like C6 (white/no special fill, centered, bordered)
$ws.Range("K6").Value = 2
$ws.Range("C6").Copy()
$ws.Range("K6").PasteSpecial(-4122)

# --- Column width for new column K ---
$ws.Range("K1").ColumnWidth = 16.9

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 90
$ws.Range("H28").Select()

$excel.CutCopyMode = $false

Write-Host "done"
